$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44279
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 35000
$ws.Range("O2").Value = 36000
$ws.Range("P2").Value = 35667
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("R2").Value = 'Región de Arica y Parinacota'
$ws.Range("S2").Value = 1982
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44442
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 35000
$ws.Range("O3").Value = 35000
$ws.Range("P3").Value = 35000
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("R3").Value = 'Perú'
$ws.Range("S3").Value = 1944
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44431
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 35000
$ws.Range("O4").Value = 35000
$ws.Range("P4").Value = 35000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Región de Arica y Parinacota'
$ws.Range("S4").Value = 1944
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44405
$ws.Range("M5").Value = 10
$ws.Range("N5").Value = 35000
$ws.Range("O5").Value = 35000
$ws.Range("P5").Value = 35000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Región de Arica y Parinacota'
$ws.Range("S5").Value = 1944
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44379
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 30000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 30000
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("R6").Value = 'Región de Arica y Parinacota'
$ws.Range("S6").Value = 1667
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44364
$ws.Range("M7").Value = 90
$ws.Range("N7").Value = 1700
$ws.Range("O7").Value = 1700
$ws.Range("P7").Value = 1700
$ws.Range("Q7").Value = '$/kilo'
$ws.Range("R7").Value = 'Región de Arica y Parinacota'
$ws.Range("S7").Value = 1700
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 44418
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 35000
$ws.Range("O8").Value = 35000
$ws.Range("P8").Value = 35000
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("R8").Value = 'Región de Arica y Parinacota'
$ws.Range("S8").Value = 1944
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44424
$ws.Range("M9").Value = 15
$ws.Range("N9").Value = 35000
$ws.Range("O9").Value = 35000
$ws.Range("P9").Value = 35000
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = 'Región de Arica y Parinacota'
$ws.Range("S9").Value = 1944
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44449
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 38000
$ws.Range("O10").Value = 38000
$ws.Range("P10").Value = 38000
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("R10").Value = 'Región de Arica y Parinacota'
$ws.Range("S10").Value = 2111
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44264
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 40000
$ws.Range("O11").Value = 40000
$ws.Range("P11").Value = 40000
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("R11").Value = 'Región de Arica y Parinacota'
$ws.Range("S11").Value = 2222
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44294
$ws.Range("M12").Value = 15
$ws.Range("N12").Value = 35000
$ws.Range("O12").Value = 35000
$ws.Range("P12").Value = 35000
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("R12").Value = 'Región de Arica y Parinacota'
$ws.Range("S12").Value = 1944
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44434
$ws.Range("M13").Value = 40
$ws.Range("N13").Value = 35000
$ws.Range("O13").Value = 35000
$ws.Range("P13").Value = 35000
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("R13").Value = 'Región de Arica y Parinacota'
$ws.Range("S13").Value = 1944
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44448
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 38000
$ws.Range("O14").Value = 38000
$ws.Range("P14").Value = 38000
$ws.Range("Q14").Value = '$/caja 18 kilos'
$ws.Range("R14").Value = 'Región de Arica y Parinacota'
$ws.Range("S14").Value = 2111
$ws.Range("T14").Value = 18

# Row 15
$ws.Range("D15").Value = 44435
$ws.Range("M15").Value = 10
$ws.Range("N15").Value = 35000
$ws.Range("O15").Value = 35000
$ws.Range("P15").Value = 35000
$ws.Range("Q15").Value = '$/caja 18 kilos'
$ws.Range("R15").Value = 'Perú'
$ws.Range("S15").Value = 1944
$ws.Range("T15").Value = 18

# Row 16
$ws.Range("D16").Value = 44435
$ws.Range("M16").Value = 105
$ws.Range("N16").Value = 35000
$ws.Range("O16").Value = 35000
$ws.Range("P16").Value = 35000
$ws.Range("Q16").Value = '$/caja 18 kilos'
$ws.Range("R16").Value = 'Región de Arica y Parinacota'
$ws.Range("S16").Value = 1944
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44377
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 40000
$ws.Range("O17").Value = 40000
$ws.Range("P17").Value = 40000
$ws.Range("Q17").Value = '$/caja 18 kilos'
$ws.Range("R17").Value = 'Región de Arica y Parinacota'
$ws.Range("S17").Value = 2222
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("D18").Value = 44433
$ws.Range("M18").Value = 15
$ws.Range("N18").Value = 35000
$ws.Range("O18").Value = 35000
$ws.Range("P18").Value = 35000
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("R18").Value = 'Región de Arica y Parinacota'
$ws.Range("S18").Value = 1944
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44363
$ws.Range("M19").Value = 144
$ws.Range("N19").Value = 1700
$ws.Range("O19").Value = 1700
$ws.Range("P19").Value = 1700
$ws.Range("Q19").Value = '$/kilo'
$ws.Range("R19").Value = 'Región de Arica y Parinacota'
$ws.Range("S19").Value = 1700
$ws.Range("T19").Value = 1

# Row 20
$ws.Range("D20").Value = 44369
$ws.Range("M20").Value = 5
$ws.Range("N20").Value = 35000
$ws.Range("O20").Value = 35000
$ws.Range("P20").Value = 35000
$ws.Range("Q20").Value = '$/caja 18 kilos'
$ws.Range("R20").Value = 'Perú'
$ws.Range("S20").Value = 1944
$ws.Range("T20").Value = 18

# Row 21
$ws.Range("D21").Value = 44357
$ws.Range("M21").Value = 10
$ws.Range("N21").Value = 38000
$ws.Range("O21").Value = 38000
$ws.Range("P21").Value = 38000
$ws.Range("Q21").Value = '$/caja 18 kilos'
$ws.Range("R21").Value = 'Perú'
$ws.Range("S21").Value = 2111
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("D22").Value = 44392
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 35000
$ws.Range("O22").Value = 35000
$ws.Range("P22").Value = 35000
$ws.Range("Q22").Value = '$/caja 18 kilos'
$ws.Range("R22").Value = 'Región de Arica y Parinacota'
$ws.Range("S22").Value = 1944
$ws.Range("T22").Value = 18

# Row 23
$ws.Range("D23").Value = 44438
$ws.Range("M23").Value = 25
$ws.Range("N23").Value = 35000
$ws.Range("O23").Value = 35000
$ws.Range("P23").Value = 35000
$ws.Range("Q23").Value = '$/caja 18 kilos'
$ws.Range("R23").Value = 'Región de Arica y Parinacota'
$ws.Range("S23").Value = 1944
$ws.Range("T23").Value = 18

# Row 24
$ws.Range("D24").Value = 44432
$ws.Range("M24").Value = 10
$ws.Range("N24").Value = 35000
$ws.Range("O24").Value = 35000
$ws.Range("P24").Value = 35000
$ws.Range("Q24").Value = '$/caja 18 kilos'
$ws.Range("R24").Value = 'Perú'
$ws.Range("S24").Value = 1944
$ws.Range("T24").Value = 18
